# Fills in the weekly teaching-routine sheet with Counselling Hour (CnH)
# and Office Hour (OH) entries for the Undergraduate Program table
# (rows 14-18, columns B/D/F/H/J/L) so that the "CnH*"/"OH*" COUNTIF
# based summary formulas in row 27 (Counselling Hour / Office Hour /
# Total Weekly Hours) recalculate correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SAT (row 14): CnH in the 1st, 4th and 5th time slots
$ws.Range("B14").Value = "CnH"
$ws.Range("J14").Value = "CnH"
$ws.Range("L14").Value = "CnH"

# SUN (row 15): CnH in the 2nd, 3rd and 4th time slots
$ws.Range("D15").Value = "CnH"
$ws.Range("H15").Value = "CnH"
$ws.Range("J15").Value = "CnH"

# MON (row 16): OH in the 2nd, 3rd and 4th time slots
$ws.Range("D16").Value = "OH"
$ws.Range("F16").Value = "OH"
$ws.Range("H16").Value = "OH"

# TUE (row 17): CnH in the 1st, 3rd and 4th time slots
$ws.Range("B17").Value = "CnH"
$ws.Range("F17").Value = "CnH"
$ws.Range("H17").Value = "CnH"

# WED (row 18): CnH in the 2nd, 4th and 5th time slots
$ws.Range("D18").Value = "CnH"
$ws.Range("H18").Value = "CnH"
$ws.Range("J18").Value = "CnH"

# Row 27 (Class/Counselling/Office/Total hour summary) recalculates
# automatically from the COUNTIF formulas already present in the sheet.
